# Atualização de bases das ligas, do dia: 2024-02-11 às 04:25
# Rotates the match-record data (columns B, F:AC) across rows 236-239,
# while keeping the row index (column A) and shared Div/Date columns
# (C, D, E) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for each affected row, keyed by row number.
# Column order: B, F, G, H, I, J, K, L, M, N, O, P, Q, R, S, T, U, V, W, X, Y, Z, AA, AB, AC
$rows = @{
    236 = @(7013409, "Nacional De Football", "Torque", 1, 1, "D",
            1.666, 3.9, 4.5, 1.615, 4, 4.75, -0.75,
            1.8, 2.05, 2.75, 1.95, 1.9,
            -1, 3, -1, -1, 1.05, -1, 0.8999999999999999)
    237 = @(7013885, "La Luz", "Atletico Fenix Montevideo", 0, 2, "A",
            3, 3, 2.4, 2.9, 2.75, 2.6, 0,
            2.025, 1.825, 2, 2.025, 1.825,
            -1, -1, 1.6, -1, 0.825, 0, 0)
    238 = @(7013886, "Racing Club de Montevideo", "Cerro", 0, 1, "A",
            2.25, 3.1, 3.25, 2.25, 2.875, 3.5, -0.25,
            1.95, 1.9, 2, 1.925, 1.925,
            -1, -1, 2.5, -1, 0.8999999999999999, -1, 0.925)
    239 = @(7013702, "Defensor Sporting", "Danubio", 0, 2, "A",
            1.8, 3.6, 4.2, 1.8, 3.6, 4.2, -0.75,
            2.05, 1.8, 2.25, 1.85, 2,
            -1, -1, 3.2, -1, 0.8, -0.5, 0.5)
}

$cols = @("B", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V", "W", "X", "Y", "Z", "AA", "AB", "AC")

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
